# Update cryptos list (Price + Volume(1h) columns) per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.064.93'
$ws.Range('E2').Value = '  -4.03%  '
$ws.Range('D3').Value = '3.653.75'
$ws.Range('E3').Value = '  -5.24%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '''589.71'
$ws.Range('E5').Value = '  -1.50%  '
$ws.Range('D6').Value = '''177.52'
$ws.Range('E6').Value = '  +4.98%  '
$ws.Range('D7').Value = '3.646.83'
$ws.Range('E7').Value = '  -5.28%  '
$ws.Range('E8').Value = '  -5.94%  '
$ws.Range('E9').Value = '  -0.08%  '
$ws.Range('D10').Value = '''0.709'
$ws.Range('E10').Value = '  -4.70%  '
$ws.Range('E11').Value = '  -8.75%  '
$ws.Range('D12').Value = '''55.05'
$ws.Range('E12').Value = '  +3.23%  '
$ws.Range('E13').Value = '  -9.81%  '
$ws.Range('D14').Value = '''10.59'
$ws.Range('E14').Value = '  -6.70%  '
$ws.Range('D15').Value = '4.231.27'
$ws.Range('E15').Value = '  -5.15%  '
$ws.Range('D16').Value = '3.649.84'
$ws.Range('E16').Value = '  -5.37%  '
$ws.Range('D17').Value = '''19.21'
$ws.Range('E17').Value = '  -9.55%  '
$ws.Range('E18').Value = '  -2.24%  '
$ws.Range('E19').Value = '  -7.21%  '
$ws.Range('D20').Value = '''12.67'
$ws.Range('E20').Value = '  -8.87%  '
$ws.Range('D21').Value = '67.888.65'
$ws.Range('E21').Value = '  -3.93%  '
$ws.Range('D22').Value = '''408.05'
$ws.Range('E22').Value = '  -7.12%  '
$ws.Range('D23').Value = '''4.54'
$ws.Range('E23').Value = '  -4.86%  '
$ws.Range('E24').Value = '  -6.70%  '
$ws.Range('E25').Value = '  -9.37%  '
$ws.Range('D26').Value = '''12.56'
$ws.Range('E26').Value = '  -9.32%  '
$ws.Range('D27').Value = '''10.74'
$ws.Range('E27').Value = '  -7.25%  '
$ws.Range('E28').Value = '  -2.96%  '
$ws.Range('E29').Value = '  +1.59%  '
$ws.Range('D30').Value = '''9.50'
$ws.Range('E30').Value = '  -9.55%  '
$ws.Range('E31').Value = '  -7.13%  '
$ws.Range('D32').Value = '''7.20'
$ws.Range('E32').Value = '  -15.38%  '
$ws.Range('E33').Value = '  -9.18%  '
$ws.Range('E34').Value = '  -6.82%  '
$ws.Range('D35').Value = '''64.22'
$ws.Range('E35').Value = '  -6.66%  '
$ws.Range('D36').Value = '''599.29'
$ws.Range('E36').Value = '  -6.27%  '
$ws.Range('D37').Value = '''42.60'
$ws.Range('E37').Value = '  -12.02%  '
$ws.Range('D38').Value = '0.0₃0881'
$ws.Range('E38').Value = '  -10.54%  '
$ws.Range('E39').Value = '  -0.06%  '
$ws.Range('D40').Value = '''0.394'
$ws.Range('E40').Value = '  -9.20%  '
$ws.Range('D41').Value = '''1.00'
$ws.Range('E41').Value = '  +0.14%  '
$ws.Range('D42').Value = '''0.136'
$ws.Range('E42').Value = '  -6.48%  '
$ws.Range('E43').Value = '  -7.37%  '
$ws.Range('D44').Value = '''2.66'
$ws.Range('E44').Value = '  -7.43%  '
$ws.Range('D45').Value = '''0.0436'
$ws.Range('E45').Value = '  -7.10%  '
$ws.Range('D46').Value = '''2.79'
$ws.Range('E46').Value = '  -11.40%  '
$ws.Range('D47').Value = '''2.71'
$ws.Range('E47').Value = '  -6.45%  '
$ws.Range('E48').Value = '  -6.82%  '
$ws.Range('D49').Value = '''8.93'
$ws.Range('E49').Value = '  -11.53%  '
$ws.Range('D50').Value = '2.706.90'
$ws.Range('E50').Value = '  -7.14%  '
$ws.Range('E51').Value = '  -6.85%  '
